$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activities Percentage")
$r = $ws.Range("A20")
$r.Interior.TintAndShade = 0.8
$r.Interior.ThemeColor = 10
Write-Host "ThemeColor after set:" $r.Interior.ThemeColor
Write-Host "TintAndShade after set:" $r.Interior.TintAndShade
Write-Host "Color:" $r.Interior.Color
